# Update forecast error table values (ifo GDP component analysis preprocessing)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q1)
$ws.Range("B2").Value = 0.02037664140334729
$ws.Range("C2").Value = 0.4125673874451535
$ws.Range("D2").Value = 0.4288236258248655
$ws.Range("E2").Value = 0.6548462612131687
$ws.Range("F2").Value = 0.6609149592609375
$ws.Range("G2").Value = 52

# Row 3 (Q2)
$ws.Range("B3").Value = 0.4198709658825853
$ws.Range("C3").Value = 0.6887531942785867
$ws.Range("D3").Value = 1.402911000708223
$ws.Range("E3").Value = 1.184445440156795
$ws.Range("F3").Value = 1.118548952961554
$ws.Range("G3").Value = 51

# Row 4 (Q3)
$ws.Range("B4").Value = 0.5914175045278451
$ws.Range("C4").Value = 1.023480694380445
$ws.Range("D4").Value = 4.028787767966156
$ws.Range("E4").Value = 2.007184039386064
$ws.Range("F4").Value = 1.937548709935223
$ws.Range("G4").Value = 50

# Row 5 (Q4)
$ws.Range("B5").Value = 0.4886808719290308
$ws.Range("C5").Value = 1.0759500315233
$ws.Range("D5").Value = 4.644551171805161
$ws.Range("E5").Value = 2.155122078167536
$ws.Range("F5").Value = 2.120737718949353
$ws.Range("G5").Value = 49

# Row 6 (Q5)
$ws.Range("B6").Value = 0.4109263009013963
$ws.Range("C6").Value = 0.9592852530966729
$ws.Range("D6").Value = 4.304132647662862
$ws.Range("E6").Value = 2.074640365861723
$ws.Range("F6").Value = 2.055056383234047
$ws.Range("G6").Value = 48

# Row 7 (Q6)
$ws.Range("B7").Value = 0.4014302792999628
$ws.Range("C7").Value = 1.035766140078607
$ws.Range("D7").Value = 5.219666336070993
$ws.Range("E7").Value = 2.284658910225111
$ws.Range("F7").Value = 2.28101933610867
$ws.Range("G7").Value = 36

# Row 8 (Q7)
$ws.Range("B8").Value = 0.3912910724344761
$ws.Range("C8").Value = 1.034972397376921
$ws.Range("D8").Value = 5.326670609579055
$ws.Range("E8").Value = 2.307958103947958
$ws.Range("F8").Value = 2.307753342904947
$ws.Range("G8").Value = 35

# Row 9 (Q8)
$ws.Range("B9").Value = 0.2626343886760503
$ws.Range("C9").Value = 1.460263562571626
$ws.Range("D9").Value = 9.616065425151794
$ws.Range("E9").Value = 3.100978140063518
$ws.Range("F9").Value = 3.179415362967537
$ws.Range("G9").Value = 18

# Row 10 (Q9)
$ws.Range("B10").Value = -0.5588777053075298
$ws.Range("C10").Value = 1.182169283809976
$ws.Range("D10").Value = 6.428227015842427
$ws.Range("E10").Value = 2.535394844169725
$ws.Range("F10").Value = 2.593736879289778
$ws.Range("G10").Value = 11
